$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -22.12810000000001
$ws.Range("C18").Value = -11.60529999999999
$ws.Range("A21").Value = -20.30029999999999
$ws.Range("A23").Value = -20.33199999999998
$ws.Range("B24").Value = 6.080699999999995
$ws.Range("A25").Value = -21.86969999999999
$ws.Range("B28").Value = 5.804600000000002
$ws.Range("B36").Value = 9.368300000000007
$ws.Range("B45").Value = 5.038800000000003
$ws.Range("B48").Value = 5.784300000000002
$ws.Range("B49").Value = 6.427499999999994
$ws.Range("C51").Value = -11.8388
$ws.Range("B52").Value = 5.6059
$ws.Range("A53").Value = -21.91569999999999
$ws.Range("B53").Value = 5.988300000000002
$ws.Range("B54").Value = 4.883200000000004
$ws.Range("C55").Value = -13.4898
$ws.Range("A57").Value = -22.18800000000002
$ws.Range("A59").Value = -21.85199999999999
$ws.Range("C64").Value = -10.80589999999999
$ws.Range("A69").Value = -21.60569999999997
$ws.Range("B70").Value = 6.7258
$ws.Range("A79").Value = -19.8472
$ws.Range("C80").Value = -13.3775
$ws.Range("A83").Value = -21.81439999999999
$ws.Range("B86").Value = 5.092700000000001
$ws.Range("B87").Value = 5.800099999999993
$ws.Range("C92").Value = -10.648
$ws.Range("A93").Value = -21.34290000000001
$ws.Range("C94").Value = -10.574
$ws.Range("C96").Value = -10.1089
$ws.Range("B101").Value = 4.925800000000001
